$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 10 (Bolivia / BO10 / Ingavi), shifting existing
#    rows 10-233 down to 11-234.
[void]$ws.Rows.Item(10).Insert()
$ws.Range("A10").Value = "Bolivia"
$ws.Range("B10").Value = "BO10"
$ws.Range("C10").Value = "Ingavi"

# 2) The Bolivia / "Cuerpo de agua" row (old row 86) is now row 87. Its PCODE
#    value needs to change from BO00 to BO11.
$ws.Range("B87").Value = "BO11"

# 3) Insert a new row at position 157 (Uruguay / UY.CH / Chuy), shifting rows
#    157-234 down to 158-235.
[void]$ws.Rows.Item(157).Insert()
$ws.Range("A157").Value = "Uruguay"
$ws.Range("B157").Value = "UY.CH"
$ws.Range("C157").Value = "Chuy"
